$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-28 Friday" "2025-11-29 Saturday"

Replace-Text "804÷8=" "340÷5="
Replace-Text "388÷7=" "376÷5="
Replace-Text "196÷8=" "768÷8="
Replace-Text "908÷6=" "142÷3="
Replace-Text "162÷8=" "787÷7="
Replace-Text "854÷3=" "255÷3="
Replace-Text "277÷4=" "109÷2="
Replace-Text "387÷3=" "475÷6="
Replace-Text "878÷7=" "812÷2="
Replace-Text "444÷3=" "586÷5="
Replace-Text "503÷6=" "439÷5="
Replace-Text "127÷3=" "272÷4="
Replace-Text "923÷8=" "172÷8="
Replace-Text "120÷8=" "484÷7="
Replace-Text "842÷6=" "348÷6="
Replace-Text "767÷2=" "845÷5="
Replace-Text "710÷3=" "873÷3="
Replace-Text "539÷9=" "933÷8="
Replace-Text "477÷6=" "471÷3="
Replace-Text "430÷3=" "763÷3="
Replace-Text "839÷4=" "322÷8="
Replace-Text "892÷8=" "742÷9="
Replace-Text "343÷9=" "157÷6="
Replace-Text "398÷7=" "151÷8="
Replace-Text "519÷7=" "927÷9="
